$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.157.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.562.29'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.20%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.79'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.561.12'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.28%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '8.01'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.166.17'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.17%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.11'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.561.37'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.261.92'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.40'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +10.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.20'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.05'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.44%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.93'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.703.28'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.04%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000117'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +8.21%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.89%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.12'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.53'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.555.98'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.153'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.04%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.90'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.86%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '170.10'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.20'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.896'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.94'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.16'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.87%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.96'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.21'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.86%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.33%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.51'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +16.48%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.14'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.81%  '
